# Workbook / worksheet handles
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at position 203.
#    This shifts the existing rows 203:295 down to 204:296 (values + formatting),
#    growing the sheet's dimension from A1:R295 to A1:R296.
$ws.Rows.Item(203).Insert()

# 2) Populate the newly-inserted (blank) row 203 with the data that used to live
#    in row 202 before the edit (i.e. the row that is being "pushed down").
$ws.Range("A203").Value = 3
$ws.Range("B203").Value = "Femacal de La Calera"
$ws.Range("C203").Value = "Coquimbo"
$ws.Range("D203").Value = 44698
$ws.Range("E203").Value = 5
$ws.Range("F203").Value = 100112013
$ws.Range("G203").Value = "Alcachofa"
$ws.Range("H203").Value = "Argentina(o)"
$ws.Range("I203").Value = "Primera"
$ws.Range("J203").Value = 145
$ws.Range("K203").Value = 16500
$ws.Range("L203").Value = 17000
$ws.Range("M203").Value = 16741
$ws.Range("N203").Value = "$/caja 50 unidades"
$ws.Range("O203").Value = "Provincia de Limarí"
$ws.Range("P203").Value = 335
$ws.Range("Q203").Value = 50
$ws.Range("R203").Value = "Hortaliza"

# 3) Update row 202 itself with the new data point values.
$ws.Range("D202").Value = 44704
$ws.Range("J202").Value = 50
$ws.Range("K202").Value = 17000
$ws.Range("M202").Value = 17000
$ws.Range("P202").Value = 340
